# Update workbook per commit: "add: nueva versión de la guía 0.1.3"
$wb = $excel.ActiveWorkbook

# 1. Bump the Version value on the Metadata sheet from 0.1.2 to 0.1.3
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B3").Value = "0.1.3"

# 2. Trim trailing whitespace from a handful of Display values on the Concepts sheet
$conceptsSheet = $wb.Worksheets.Item("Concepts")

$conceptsSheet.Range("C5").Value = "SEREMI De Arica y Parinacota Clínica Hebe"
$conceptsSheet.Range("C6").Value = "SEREMI De Arica y Parinacota Clínica San Agustín"
$conceptsSheet.Range("C7").Value = "SEREMI De Arica y Parinacota Clínica San José"
$conceptsSheet.Range("C9").Value = "SEREMI De Arica y Parinacota Complejo Penitenciario"
